$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert row 95 with a new health facility entry (Pebane / PS Cutal)
$ws.Rows(95).Insert()
$ws.Range("B95").Value = "ZAMBEZIA"
$ws.Range("C95").Value = "Pebane"
$ws.Range("D95").Value = "PS Cutal"

# Insert row 96 with a new health facility entry (Namacurra / PS Naciaia)
$ws.Rows(96).Insert()
$ws.Range("B96").Value = "ZAMBEZIA"
$ws.Range("C96").Value = "Namacurra"
$ws.Range("D96").Value = "PS Naciaia"

# Restore the view/selection state recorded for the edited sheet
$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D84").Select()
